# Apply updates described in the commit "update to published CDA FHIR
# logical model with patches #241":
#  - Metadata sheet: Version, Date, Contact values updated
#  - Elements sheet: Qualifier binding value set URL updated
#  - Elements sheet: column Z (Binding Value Set) width widened

$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B3").Value = "2.0.0-sd-202406-matchbox-patch"
$wsMeta.Range("B8").Value = "2024-06-19T17:47:42+02:00"
$wsMeta.Range("B10").Value = "HL7 International - Structured Documents (http://www.hl7.org/Special/committees/structure, structdog@lists.HL7.org)"

$wsElements = $wb.Worksheets.Item("Elements")
$wsElements.Range("Z9").Value = "http://hl7.org/cda/stds/core/ValueSet/CDAEntityNamePartQualifier"
# Target stored column width is 61.7578125 characters; the COM ColumnWidth
# setter here rounds pixel width to a 1/6-character grid, so feed it the
# character width (minus the 5px/6 padding offset) that lands on the closest
# achievable grid point to the target.
$wsElements.Columns.Item(26).ColumnWidth = 60.924479166666664
